$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in row 9 with the new "h1-dbt" heading-size row, extending the
# existing table pattern (rows 2-8) down by one more record.
$ws.Range("A9").Value = "h1-dbt"
$ws.Range("B9").Value = 1.75
$ws.Range("C9").Value = 1
$ws.Range("D9").Value = 1.35

$ws.Range("E9").Formula = '=IF(($E$1/100*B9/16)+C9<D9,D9,IF(($E$1/100*B9/16)+C9>J9,J9,($E$1/100*B9/16)+C9))'
$ws.Range("F9").Formula = '=IF(($F$1/100*B9/16)+C9<D9,D9,IF(($F$1/100*B9/16)+C9>J9,J9,($F$1/100*B9/16)+C9))'
$ws.Range("G9").Formula = '=IF(($G$1/100*B9/16)+C9<D9,D9,IF(($G$1/100*B9/16)+C9>J9,J9,($G$1/100*B9/16)+C9))'
$ws.Range("H9").Formula = '=IF(($H$1/100*B9/16)+C9<D9,D9,IF(($H$1/100*B9/16)+C9>J9,J9,($H$1/100*B9/16)+C9))'
$ws.Range("I9").Formula = '=IF(($I$1/100*B9/16)+C9<D9,D9,IF(($I$1/100*B9/16)+C9>J9,J9,($I$1/100*B9/16)+C9))'

$ws.Range("J9").Value = 3
$ws.Range("K9").Value = 2.5

$ws.Range("L9").Formula = '=E9*16'
$ws.Range("M9").Formula = '=F9*16'
$ws.Range("N9").Formula = '=G9*16'
$ws.Range("O9").Formula = '=H9*16'
$ws.Range("P9").Formula = '=I9*16'

$ws.Range("S9").Formula = '="font-size: clamp("&D9&"rem, "&B9&"vw + "&C9&"rem, "&J9&"rem);"'

# Entering formulas into cells that previously matched row 8's shared-formula
# pattern re-paints L9:P9 with row 8's fill color. Restore the original
# banding (matching rows 10+) by re-pasting the untouched formatting from
# the row below, which keeps the same style without minting a new one.
$ws.Range("L10:P10").Copy()
$ws.Range("L9:P9").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Update the active selection to A10, matching the post-edit sheet view.
$ws.Range("A10").Select()
